$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.32"
$ws.Range("E2").Value = "'1.08%"
$ws.Range("D3").Value = "'36.31"
$ws.Range("E3").Value = "'1.17%"
$ws.Range("E4").Value = "'0.63%"
$ws.Range("D5").Value = "'0.08109"
$ws.Range("E5").Value = "'0.16%"
$ws.Range("E6").Value = "'7.83%"
$ws.Range("D7").Value = "'4.156"
$ws.Range("E7").Value = "'0.35%"
$ws.Range("E8").Value = "'0.05%"
$ws.Range("D9").Value = "'0.9302"
$ws.Range("E9").Value = "'-0.20%"
$ws.Range("D10").Value = "'0.1409"
$ws.Range("E10").Value = "'10.14%"
$ws.Range("D11").Value = "'0.1934"
$ws.Range("E11").Value = "'0.88%"
$ws.Range("D12").Value = "'0.09073"
$ws.Range("E12").Value = "'-1.31%"
$ws.Range("E13").Value = "'-1.08%"
$ws.Range("D14").Value = "'0.09919"
$ws.Range("E14").Value = "'-0.10%"
$ws.Range("D15").Value = "'0.001406"
$ws.Range("E15").Value = "'-0.67%"
$ws.Range("D16").Value = "'0.006315"
$ws.Range("E16").Value = "'-4.81%"
$ws.Range("E17").Value = "'6.25%"
$ws.Range("D18").Value = "'3.336"
$ws.Range("E18").Value = "'4.25%"
$ws.Range("D19").Value = "'0.3446"
$ws.Range("E19").Value = "'0.00%"
$ws.Range("D20").Value = "'0.1308"
$ws.Range("E20").Value = "'0.27%"
$ws.Range("E21").Value = "'-7.31%"
$ws.Range("D22").Value = "'0.2341"
$ws.Range("E22").Value = "'-7.58%"
$ws.Range("D23").Value = "'0.04376"
$ws.Range("E23").Value = "'-0.87%"
$ws.Range("E24").Value = "'-0.23%"
$ws.Range("E25").Value = "'4.08%"
$ws.Range("D27").Value = "'0.0001299"
$ws.Range("E27").Value = "'-0.21%"
$ws.Range("E39").Value = "'2.47%"
$ws.Range("D40").Value = "'0.05163"
$ws.Range("E40").Value = "'-0.62%"
$ws.Range("D41").Value = "'0.007499"
$ws.Range("E41").Value = "'-0.49%"
$ws.Range("D42").Value = "'0.01003"
$ws.Range("E42").Value = "'-0.66%"
$ws.Range("E43").Value = "'-0.03%"
$ws.Range("D44").Value = "'0.002169"
$ws.Range("E44").Value = "'3.12%"
$ws.Range("D45").Value = "'0.009977"
$ws.Range("E45").Value = "'-6.68%"
$ws.Range("D46").Value = "'0.00006284"
$ws.Range("E46").Value = "'-1.13%"
$ws.Range("E47").Value = "'-0.13%"
$ws.Range("D48").Value = "'64.85"
$ws.Range("E48").Value = "'-0.16%"
$ws.Range("E49").Value = "'-22.03%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.13%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.13%"
